# Apply refreshed cryptocurrency price/volume data (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.603.48"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "3.281.16"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "535.81"
$ws.Range("E5").Value = "  +5.10%  "
$ws.Range("D6").Value = "173.95"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("D7").Value = "0.597"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("D8").Value = "3.276.89"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.611"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "'53.80"
$ws.Range("E11").Value = "  -5.36%  "
$ws.Range("D12").Value = "0.135"
$ws.Range("E12").Value = "  +4.69%  "
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").Value = "9.16"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "3.812.16"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "3.281.46"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "17.39"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "63.585.92"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "11.18"
$ws.Range("E20").Value = "  +3.39%  "
$ws.Range("D21").Value = "0.968"
$ws.Range("E21").Value = "  +3.50%  "
$ws.Range("D22").Value = "370.41"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").Value = "11.43"
$ws.Range("E23").Value = "  +3.47%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "4.12"
$ws.Range("E24").Value = "  +9.66%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "3.78"
$ws.Range("E25").Value = "  +4.91%  "
$ws.Range("D26").Value = "81.37"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("D27").Value = "6.17"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("D29").Value = "11.36"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").Value = "'8.30"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").Value = "28.72"
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("D32").Value = "642.16"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").Value = "6.49"
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("D34").Value = "11.29"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("E35").Value = "  +4.97%  "
$ws.Range("D36").Value = "56.95"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'36.90"
$ws.Range("E38").Value = "  +3.65%  "
$ws.Range("D39").Value = "0.384"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("E40").Value = "  +14.93%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "2.64"
$ws.Range("E42").Value = "  +9.54%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.125"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "2.925.11"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("E45").Value = "  +6.35%  "
$ws.Range("D46").Value = "'2.70"
$ws.Range("E46").Value = "  +6.12%  "
$ws.Range("E47").Value = "  +5.88%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "2.64"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "3.09"
$ws.Range("E49").Value = "  +5.78%  "
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").Value = "'135.50"
$ws.Range("E51").Value = "  +5.24%  "
